$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.785.06"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.137.01"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'586.76"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'146.28"
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.135.13"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "'0.458"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").Value = "'36.85"
$ws.Range("E14").Value = "  +2.46%  "
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "3.651.89"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "63.584.33"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "3.134.21"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'7.08"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "'464.33"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "'14.26"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'7.42"
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("D25").Value = "'81.12"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'9.27"
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "'7.00"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'26.97"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  -5.91%  "
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").Value = "'6.01"
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "'51.27"
$ws.Range("D41").Value = "'440.42"
$ws.Range("E41").Value = "  -2.76%  "
$ws.Range("D42").Value = "'8.80"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.918.02"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0370"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "'0.278"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").Value = "'36.93"
$ws.Range("E47").Value = "  +3.53%  "
$ws.Range("D48").Value = "'126.76"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("D51").Value = "'24.13"
$ws.Range("E51").Value = "  -3.89%  "
